$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.366.11'
$ws.Range("E2").Value = '  -1.38%  '
$ws.Range("D3").Value = '3.421.41'
$ws.Range("E3").Value = '  -2.12%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.45'
$ws.Range("E5").Value = '  -2.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.53'
$ws.Range("E6").Value = '  -3.67%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -1.54%  '
$ws.Range("E9").Value = '  +2.16%  '
$ws.Range("E10").Value = '  -1.22%  '
$ws.Range("E11").Value = '  -1.23%  '
$ws.Range("D12").Value = '4.003.33'
$ws.Range("E12").Value = '  -2.14%  '
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("E14").Value = '  -2.71%  '
$ws.Range("D15").Value = '3.419.98'
$ws.Range("E15").Value = '  -2.19%  '
$ws.Range("D16").Value = '63.372.23'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.05'
$ws.Range("E17").Value = '  -2.16%  '
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("E19").Value = '  -1.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.18'
$ws.Range("E20").Value = '  -2.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '382.96'
$ws.Range("E21").Value = '  -2.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.561'
$ws.Range("E22").Value = '  -1.89%  '
$ws.Range("D23").Value = '3.559.71'
$ws.Range("E23").Value = '  -2.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.60'
$ws.Range("E24").Value = '  -1.29%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").Value = '  -4.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.993'
$ws.Range("E27").Value = '  -0.56%  '
$ws.Range("E28").Value = '  -3.33%  '
$ws.Range("E29").Value = '  -4.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.88'
$ws.Range("E30").Value = '  -3.82%  '
$ws.Range("E31").Value = '  -0.63%  '
$ws.Range("E32").Value = '  -4.50%  '
$ws.Range("D33").Value = '3.451.90'
$ws.Range("E33").Value = '  -1.81%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '22.61'
$ws.Range("E35").Value = '  -3.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.16'
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.73'
$ws.Range("E37").Value = '  -2.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '163.65'
$ws.Range("E38").Value = '  -1.77%  '
$ws.Range("E39").Value = '  -2.60%  '
$ws.Range("E40").Value = '  -2.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.784'
$ws.Range("E41").Value = '  -3.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.20'
$ws.Range("E43").Value = '  -1.15%  '
$ws.Range("E44").Value = '  -2.11%  '
$ws.Range("E45").Value = '  -3.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.35'
$ws.Range("E46").Value = '  -6.53%  '
$ws.Range("E47").Value = '  -6.14%  '
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.880'
$ws.Range("E49").Value = '  -1.16%  '
$ws.Range("D50").Value = '2.257.73'
$ws.Range("E50").Value = '  -4.77%  '
$ws.Range("E51").Value = '  -2.62%  '
